$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each group lists worksheet row numbers whose B/E/F/G values are cyclically
# rotated: row[i] takes the B/E/F/G values that row[i+1] currently holds
# (wrapping around to the first row of the group). Columns A (index),
# C (item name) and D (purchase rate) stay attached to their own row.
$groups = @(
    @(161,162,163),
    @(183,184),
    @(264,265),
    @(351,352),
    @(355,356),
    @(375,376),
    @(379,380),
    @(400,401),
    @(579,580),
    @(583,584),
    @(590,591),
    @(593,594),
    @(601,602),
    @(687,688),
    @(709,710),
    @(715,716)
)

foreach ($grp in $groups) {
    $n = $grp.Length

    # Snapshot current values for the columns that rotate.
    $bVals = @()
    $eVals = @()
    $fVals = @()
    $gVals = @()
    for ($i = 0; $i -lt $n; $i++) {
        $r = $grp[$i]
        $bVals += $ws.Cells.Item($r, 2).Value2
        $eVals += $ws.Cells.Item($r, 5).Value2
        $fVals += $ws.Cells.Item($r, 6).Value2
        $gVals += $ws.Cells.Item($r, 7).Value2
    }

    # Write back: row i gets the values that belonged to row (i+1) mod n.
    for ($i = 0; $i -lt $n; $i++) {
        $r = $grp[$i]
        $src = ($i + 1) % $n
        $ws.Cells.Item($r, 2).Value2 = $bVals[$src]
        $ws.Cells.Item($r, 5).Value2 = $eVals[$src]
        $ws.Cells.Item($r, 6).Value2 = $fVals[$src]
        $ws.Cells.Item($r, 7).Value2 = $gVals[$src]
    }
}
